$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value. Values for the "Price" (D) column
# that look like plain decimal numbers are prefixed with a leading
# apostrophe so Excel stores them as literal text (preserving trailing
# zeros / exact digit formatting like "1.00" or "38.20") instead of
# silently converting them to numeric values. The "Volume(1h)" (E) column
# values are percentage strings with surrounding spaces, which Excel
# always keeps as text already.
$updates = @{
    'D2' = '66.813.61'
    'E2' = '  -4.99%  '
    'D3' = '3.205.00'
    'E3' = '  -8.76%  '
    'D4' = '''1.00'
    'E4' = '  -0.09%  '
    'D5' = '''585.91'
    'E5' = '  -3.14%  '
    'D6' = '''148.27'
    'E6' = '  -15.13%  '
    'E7' = '  -0.01%  '
    'D8' = '3.199.67'
    'E8' = '  -8.67%  '
    'D9' = '''0.534'
    'E9' = '  -12.69%  '
    'D10' = '''0.168'
    'E10' = '  -13.03%  '
    'D11' = '''6.28'
    'E11' = '  -12.58%  '
    'D12' = '''0.477'
    'E12' = '  -18.14%  '
    'D13' = '''38.20'
    'E13' = '  -17.29%  '
    'D14' = '''0.0000238'
    'E14' = '  -13.38%  '
    'D15' = '3.735.98'
    'E15' = '  -8.28%  '
    'D16' = '66.784.64'
    'E16' = '  -5.27%  '
    'D17' = '3.204.84'
    'E17' = '  -9.07%  '
    'E18' = '  -5.17%  '
    'D19' = '''516.26'
    'E19' = '  -15.31%  '
    'D20' = '''6.85'
    'E20' = '  -17.28%  '
    'D21' = '''14.29'
    'E21' = '  -17.91%  '
    'D22' = '''0.740'
    'E22' = '  -15.55%  '
    'D23' = '''7.73'
    'E23' = '  -14.10%  '
    'D24' = '''83.59'
    'E24' = '  -15.06%  '
    'D25' = '''13.17'
    'E25' = '  -15.34%  '
    'D27' = '''3.14'
    'E27' = '  -15.59%  '
    'D28' = '''2.09'
    'E28' = '  -18.05%  '
    'D29' = '''7.79'
    'E29' = '  -13.34%  '
    'D30' = '''28.36'
    'E30' = '  -16.08%  '
    'E31' = '  -12.80%  '
    'D32' = '''1.14'
    'E32' = '  -10.71%  '
    'D33' = '''528.07'
    'E33' = '  -16.76%  '
    'D34' = '''5.60'
    'E34' = '  -17.91%  '
    'D35' = '''1.00'
    'E35' = '  +0.21%  '
    'D36' = '''6.29'
    'E36' = '  -21.70%  '
    'D37' = '''53.35'
    'E37' = '  -5.88%  '
    'D38' = '''0.0423'
    'E38' = '  -10.69%  '
    'D39' = '''0.0846'
    'E39' = '  -14.78%  '
    'D40' = '''9.05'
    'E40' = '  -15.80%  '
    'D41' = '''0.122'
    'E41' = '  -14.11%  '
    'D42' = '''2.69'
    'E42' = '  -24.70%  '
    'D43' = '2.853.84'
    'E43' = '  -15.39%  '
    'D44' = '0.0₃0573'
    'E44' = '  -22.60%  '
    'D45' = '''0.254'
    'E45' = '  -17.42%  '
    'D46' = '''2.37'
    'E46' = '  -18.03%  '
    'E47' = '  -0.09%  '
    'D48' = '''2.10'
    'E48' = '  -17.52%  '
    'D49' = '''25.39'
    'E49' = '  -21.04%  '
    'D50' = '''0.112'
    'E50' = '  -13.40%  '
    'D51' = '''120.35'
    'E51' = '  -9.18%  '
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
